$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Range("D18").Value = "[-, -, -, 'ELM-2NA-Tecnologia da Soldagem']"
$ws.Range("E18").Value = "-"
$ws.Range("F18").Value = "-"

# Row 19
$ws.Range("D19").Value = "[-, -, 'ELM-2NA-Tecnologia da Soldagem', -]"
$ws.Range("F19").Value = "-"

# Row 20
$ws.Range("F20").Value = "-"

# Row 21
$ws.Range("D21").Value = "[-, -, 'ELM-2NA-Tecnologia da Soldagem', -]"
$ws.Range("F21").Value = "[-, -, 'ELM-2NA-Tecnologia da Soldagem', -]"
